# Docx: corrige bug relacionado con los textboxs.
#
# 1. Remove the stray "_GoBack" bookmark that wrapped the first inline
#    textbox run (it no longer belongs there once a second textbox is
#    added further down).
# 2. Append a new paragraph ("Debajo hay otro cuadro de texto:") and a
#    new paragraph containing a second (anchored) textbox shape, whose
#    own content carries its own "_GoBack" bookmark pair.

$d = $word.ActiveDocument

$bookmark = $d.Bookmarks.Item('_GoBack')
$bookmark.Delete()

$insertXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" mc:Ignorable="w14 wp14"><w:body><w:p><w:r><w:t>Debajo hay otro cuadro de texto:</w:t></w:r></w:p><w:p><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:editId="36B11C9B"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:align>center</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>0</wp:posOffset></wp:positionV><wp:extent cx="2374265" cy="1403985"/><wp:effectExtent l="0" t="0" r="3175" b="5715"/><wp:wrapNone/><wp:docPr id="2" name="Cuadro de texto 2"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"><a:spLocks noChangeArrowheads="1"/></wps:cNvSpPr><wps:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="2374265" cy="1403985"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:solidFill><a:srgbClr val="FFFFFF"/></a:solidFill><a:ln w="9525"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></wps:spPr><wps:txbx><w:txbxContent><w:p><w:r><w:t>Blablabla1.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:sdt><w:sdtPr><w:id w:val="568603642"/><w:temporary/><w:showingPlcHdr/></w:sdtPr><w:sdtContent><w:p><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>[Escriba una cita del documento o el resumen de un punto interesante. Puede situar el cuadro de texto en cualquier lugar del documento. Use la ficha Herramientas de dibujo para cambiar el formato del cuadro de texto de la cita.]</w:t></w:r></w:p></w:sdtContent></w:sdt><w:p><w:r><w:t>Blablabla2.</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" anchor="t" anchorCtr="0"><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>40000</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>20000</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shapetype id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="_x0000_s1028" type="#_x0000_t202" style="position:absolute;margin-left:0;margin-top:0;width:186.95pt;height:110.55pt;z-index:251659264;visibility:visible;mso-wrap-style:square;mso-width-percent:400;mso-height-percent:200;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:center;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:400;mso-height-percent:200;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:top" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQA91hoqLAIAAFMEAAAOAAAAZHJzL2Uyb0RvYy54bWysVNtu2zAMfR+wfxD0vjpxk6416hRdugwD&#10;ugvQ7QMYSY6FyaImKbG7rx8lu1l2exnmB0EUqSPyHNLXN0Nn2EH5oNHWfH4240xZgVLbXc0/f9q8&#10;uOQsRLASDFpV88cV+M3q+bPr3lWqxBaNVJ4RiA1V72rexuiqogiiVR2EM3TKkrNB30Ek0+8K6aEn&#10;9M4U5Wx2UfTopfMoVAh0ejc6+SrjN40S8UPTBBWZqTnlFvPq87pNa7G6hmrnwbVaTGnAP2TRgbb0&#10;6BHqDiKwvde/QXVaeAzYxDOBXYFNo4XKNVA189kv1Ty04FSuhcgJ7khT+H+w4v3ho2da1rzkzEJH&#10;Eq33ID0yqVhUQ0RWJpJ6FyqKfXAUHYdXOJDYueDg7lF8CcziugW7U7feY98qkJTkPN0sTq6OOCGB&#10;bPt3KOk12EfMQEPju8QgccIIncR6PApEeTBBh+X5y0V5seRMkG++mJ1fXS7zG1A9XXc+xDcKO5Y2&#10;NffUARkeDvchpnSgegpJrwU0Wm60Mdnwu+3aeHYA6pZN/ib0n8KMZX3Nr5blcmTgrxCz/P0JotOR&#10;2t7oruaXxyCoEm+vrcxNGUGbcU8pGzsRmbgbWYzDdpiEm/TZonwkZj2OXU5TSZsW/TfOeurwmoev&#10;e/CKM/PWkjpX88UijUQ2FsuXJRn+1LM99YAVBFXzyNm4Xcc8Rpk3d0sqbnTmN8k9ZjKlTJ2baZ+m&#10;LI3GqZ2jfvwLVt8BAAD//wMAUEsDBBQABgAIAAAAIQD9LzLW2wAAAAUBAAAPAAAAZHJzL2Rvd25y&#10;ZXYueG1sTI/BTsMwEETvSPyDtUjcqJMUFUjjVFUE10ptkbhu420SsNchdtLw9xgucFlpNKOZt8Vm&#10;tkZMNPjOsYJ0kYAgrp3uuFHweny5ewThA7JG45gUfJGHTXl9VWCu3YX3NB1CI2IJ+xwVtCH0uZS+&#10;bsmiX7ieOHpnN1gMUQ6N1ANeYrk1MkuSlbTYcVxosaeqpfrjMFoF47HaTvsqe3+bdvp+t3pGi+ZT&#10;qdubebsGEWgOf2H4wY/oUEamkxtZe2EUxEfC743e8mH5BOKkIMvSFGRZyP/05TcAAAD//wMAUEsB&#10;Ai0AFAAGAAgAAAAhALaDOJL+AAAA4QEAABMAAAAAAAAAAAAAAAAAAAAAAFtDb250ZW50X1R5cGVz&#10;XS54bWxQSwECLQAUAAYACAAAACEAOP0h/9YAAACUAQAACwAAAAAAAAAAAAAAAAAvAQAAX3JlbHMv&#10;LnJlbHNQSwECLQAUAAYACAAAACEAPdYaKiwCAABTBAAADgAAAAAAAAAAAAAAAAAuAgAAZHJzL2Uy&#10;b0RvYy54bWxQSwECLQAUAAYACAAAACEA/S8y1tsAAAAFAQAADwAAAAAAAAAAAAAAAACGBAAAZHJz&#10;L2Rvd25yZXYueG1sUEsFBgAAAAAEAAQA8wAAAI4FAAAAAA==&#10;"><v:textbox style="mso-fit-shape-to-text:t"><w:txbxContent><w:p><w:r><w:t>Blablabla1.</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p><w:sdt><w:sdtPr><w:id w:val="568603642"/><w:temporary/><w:showingPlcHdr/></w:sdtPr><w:sdtContent><w:p><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>[Escriba una cita del documento o el resumen de un punto interesante. Puede situar el cuadro de texto en cualquier lugar del documento. Use la ficha Herramientas de dibujo para cambiar el formato del cuadro de texto de la cita.]</w:t></w:r></w:p></w:sdtContent></w:sdt><w:p><w:r><w:t>Blablabla2.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target = $d.Content
$target.Collapse(0)
$target.InsertXML($insertXml)
